{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Helper to find paragraph index whose text matches exactly\nfunction findIndex(items, text) {\n  for (let i = 0; i < items.length; i++) {\n    if (items[i].text === text) return i;\n  }\n  return -1;\n}\n\n// 1) Ativa\u00e7\u00e3o: 01/01/2021 -> Ativa\u00e7\u00e3o: 01/01/2024\nconst searchResults = body.search(\"Ativa\u00e7\u00e3o: 01/01/2021\", { matchCase: true });\nsearchResults.load(\"text\");\nawait context.sync();\nif (searchResults.items.length > 0) {\n  searchResults.items[0].insertText(\"Ativa\u00e7\u00e3o: 01/01/2024\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 2) Objetivos paragraph: replace text, then add new italic paragraph after it\nlet idx = findIndex(paragraphs.items, \"Apresentar os conceitos de Log\u00edstica, Log\u00edstica Reversa e Gest\u00e3o da Cadeia de Suprimentos. Capacitar o aluno para aplica\u00e7\u00e3o de t\u00e9cnicas e m\u00e9todos quantitativos para otimiza\u00e7\u00e3o dos problemas em Log\u00edstica e Cadeias de Suprimentos.\");\nlet target = paragraphs.items[idx];\ntarget.insertText(\"Fornecer conhecimentos que proporcionam uma vis\u00e3o hol\u00edstica a respeito da Cadeia de Suprimentos e da Log\u00edstica, apresentando m\u00e9todos e ferramentas para otimizar o desempenho das cadeias produtivas.\", Word.InsertLocation.replace);\nconst newPara1 = target.insertParagraph(\"Provide knowledge that paves a holistic view of Logistics and Supply Chain, presenting methods and tools available to optimize the performance of production chains.\", Word.InsertLocation.after);\nnewPara1.font.italic = true;\nawait context.sync();\n\n// Reload paragraphs after structural change\nconst paragraphs2 = body.paragraphs;\nparagraphs2.load(\"text\");\nawait context.sync();\n\n// 3) Docente responsible person replacement\nidx = findIndex(paragraphs2.items, \"5840560 - Marco Antonio Carvalho Pereira\");\nparagraphs2.items[idx].insertText(\"3295113 - Jos\u00e9 Eduardo Holler Branco\", Word.InsertLocation.replace);\nawait context.sync();\n\n// 4) Programa resumido paragraph: replace text, then add new italic paragraph after it\nconst paragraphs3 = body.paragraphs;\nparagraphs3.load(\"text\");\nawait context.sync();\nidx = findIndex(paragraphs3.items, \"1. Introdu\u00e7\u00e3o: 2. Gest\u00e3o estrat\u00e9gica3. Gest\u00e3o dos relacionamentos4. Gest\u00e3o global de suprimentos5. Avalia\u00e7\u00e3o de desempenho6. Mapeamento e an\u00e1lise de processos7. Gest\u00e3o de demanda8. Gest\u00e3o e coordena\u00e7\u00e3o de estoques9. Gest\u00e3o da log\u00edstica10. Log\u00edstica reversa\");\ntarget = paragraphs3.items[idx];\ntarget.insertText(\"Gerenciamento da cadeia de suprimentos e da log\u00edstica: planejamento, otimiza\u00e7\u00e3o e controle.\", Word.InsertLocation.replace);\nconst newPara2 = target.insertParagraph(\"Supply chain and logistics management: planning, optimization and control.\", Word.InsertLocation.after);\nnewPara2.font.italic = true;\nawait context.sync();\n\n// 5) Programa (full) paragraph: replace text, then add new italic paragraph after it\nconst paragraphs4 = body.paragraphs;\nparagraphs4.load(\"text\");\nawait context.sync();\nidx = findIndex(paragraphs4.items, \"1. Introdu\u00e7\u00e3o: A concorr\u00eancia entre cadeias de suprimento. Defini\u00e7\u00e3o operacional. A globaliza\u00e7\u00e3o e a gest\u00e3o de cadeia de suprimentos. Governan\u00e7a das cadeias de suprimentos2. Gest\u00e3o estrat\u00e9gica: Estrat\u00e9gia de cadeia de suprimentos. Produtos funcionais x produtos inovadores. Fluxos empurrados puxados e h\u00edbridos. Custo de transa\u00e7\u00e3o e a decis\u00e3o estrat\u00e9gica de comprar ou fazer. Padroniza\u00e7\u00e3o. Integra\u00e7\u00e3o de parceiros da cadeia de suprimento no projeto de novos produtos e processos.3. Gest\u00e3o dos relacionamentos: Confian\u00e7a entre parceiros. Negocia\u00e7\u00e3o. Gest\u00e3o do relacionamento com clientes. Segmenta\u00e7\u00e3o de produtos. Gest\u00e3o do relacionamento com fornecedores4. Gest\u00e3o global de suprimentos: Tipos de suplemento. Estrutura organizacional para suprimentos. O processo de suprimento. Coopeti\u00e7\u00e3o. \u00c9tica e responsabilidade social na gest\u00e3o global de suprimentos5. Avalia\u00e7\u00e3o de desempenho: O que \u00e9 medi\u00e7\u00e3o de desempenho? Porque medir desempenho. Caracter\u00edsticas de uma boa medida de desempenho. Alinhamento de incentivos em cadeias globais de suprimento. Tipos de contrato de relacionamento6. Mapeamento e an\u00e1lise de processos: Principais processos na cadeia de suprimento. O modelo SCOR (Supply Chain Operations Reference). An\u00e1lise e melhoramento de processos.7. Gest\u00e3o de demanda: A\u00e7\u00f5es sobre a demanda para redu\u00e7\u00e3o de variabilidade. Causas da variabilidade da demanda. Previs\u00e3o de demanda. Processo de previs\u00e3o de vendas. M\u00e9todos usados em previs\u00f5es. M\u00e9todo Delphi. Incerteza de previs\u00e3o8. Gest\u00e3o e coordena\u00e7\u00e3o de estoques: Defini\u00e7\u00e3o de estoques. Causa do surgimento dos estoques. Tipos de estoque. VMI (vendor management inventory) - estoque gerenciado pelo distribuidor. VOI (vendor owner inventory) - consigna\u00e7\u00e3o9. Gest\u00e3o da log\u00edstica: Centraliza\u00e7\u00e3o versus descentraliza\u00e7\u00e3o. Pontos de armazenagem/distribui\u00e7\u00e3o. Fun\u00e7\u00f5es dos armaz\u00e9ns. Sistemas log\u00edsticos escalonados. Localiza\u00e7\u00e3o de unidades log\u00edsticas. Gest\u00e3o de transportes na cadeia de suprimentos.10. Log\u00edstica reversa: Conceito, import\u00e2ncia, estrutura e tend\u00eancias. Sustentabilidade. Ciclo fechado. Tipos de ciclo fechado. Motiva\u00e7\u00e3o empresarial. Gerenciamento integrado de res\u00edduos. Modelos de roteiriza\u00e7\u00e3o. Programa\u00e7\u00e3o de frotas de ve\u00edculos.\");\ntarget = paragraphs4.items[idx];\ntarget.insertText(\"i) Introdu\u00e7\u00e3o \u00e0 Log\u00edstica e Cadeia de Suprimentos; ii) Planejamento da cadeia de suprimentos; iii) Planejamento do transporte; iv) Custos log\u00edsticos; v) Tipos de cargas e sistemas de armazenamento; vi) Modelos de transporte; vii) Modelos de localiza\u00e7\u00e3o; viii) Planejamento do estoque; ix) Log\u00edstica Reversa e Economia Circular; e x)  Controle da log\u00edstica e cadeia de suprimentos.\", Word.InsertLocation.replace);\nconst newPara3 = target.insertParagraph(\"i) Introduction to Logistics and Supply Chain; ii) Supply chain planning; iii) Transport planning; iv) Logistic costs; v) Types of cargos and storage systems; vi) Stock planning; vii) Transport models; viii) Location models; ix) Reverse Logistics and Circular Economy; and x) Control of logistics and supply chain.\", Word.InsertLocation.after);\nnewPara3.font.italic = true;\nawait context.sync();\n\n// 6) M\u00e9todo: text replacement\nconst searchMethod = body.search(\"Aulas expositivas te\u00f3ricas, aulas de exerc\u00edcios.\", { matchCase: true });\nsearchMethod.load(\"text\");\nawait context.sync();\nif (searchMethod.items.length > 0) {\n  searchMethod.items[0].insertText(\"Provas, trabalhos em grupo, exerc\u00edcios individuais e semin\u00e1rios.\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 7) Crit\u00e9rio: text replacement\nconst searchCriterio = body.search(\"M\u00e9dia de Provas e trabalhos (MF).\", { matchCase: true });\nsearchCriterio.load(\"text\");\nawait context.sync();\nif (searchCriterio.items.length > 0) {\n  searchCriterio.items[0].insertText(\"M\u00e9dia das atividades avaliativas.\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 8) Norma de recupera\u00e7\u00e3o: text replacement\nconst searchNorma = body.search(\"Prova de Recupera\u00e7\u00e3o (PR). A Nota final (NF) ser\u00e1 a m\u00e9dia aritm\u00e9tica entre MF e PR\", { matchCase: true });\nsearchNorma.load(\"text\");\nawait context.sync();\nif (searchNorma.items.length > 0) {\n  searchNorma.items[0].insertText(\"NF = (MF + PR)/2, onde MF \u00e9 a m\u00e9dia final da avalia\u00e7\u00e3o e PR \u00e9 uma prova de recupera\u00e7\u00e3o.\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 9) Bibliografia paragraph replacement\nconst paragraphs5 = body.paragraphs;\nparagraphs5.load(\"text\");\nawait context.sync();\nidx = findIndex(paragraphs5.items, \"CORR\u00caA, HENRIQUE LUIZ. Gest\u00e3o de rede de suprimento: integrando cadeias de suprimento no mundo globalizado. Editora Atlas, 2009CORREA, HENRIQUE LUIZ. Administra\u00e7\u00e3o de cadeias de suprimento e log\u00edstica: o essencial. Editora Atlas 2014PIRES, S\u00c9RGIO. Gest\u00e3o da cadeia de suprimentos (Supply Chain Management): conceitos, estrat\u00e9gias, pr\u00e1ticas e casos. Editora Atlas segunda edi\u00e7\u00e3o. 2009IYER, ANANTH; SESHHADRI, SHIDHAR; VASHER, ROY. A gest\u00e3o da cadeia de suprimentos da Toyota. Bookman. 2009MYERSON, PAUL. Lean Supply Chain and logistics management. McGrawHill. 2012\");\nparagraphs5.items[idx].insertText(\"BOWERSOX, D. J.; CLOSS, D. J.; COOPER, M. B.; BOWERSOX, J. C. Gest\u00e3o Log\u00edstica da Cadeia de Suprimentos. 4. ed. AMGH, 2013. 472 p.BARTHOLOMEU, D. B.; CAIXETA FILHO, J. V. Log\u00edstica Ambiental de Res\u00edduos S\u00f3lidos. S\u00e3o Paulo: Atlas, 2011, 249 p.CHOPRA, S.; MEINDL, P. Gest\u00e3o da cadeia de suprimentos: estrat\u00e9gia, planejamento e Opera\u00e7\u00f5es. 6. ed. Pearson, 2015. 544 p.CAIXETA FILHO, J. V.; MARTINS, R. S. (org.). Gest\u00e3o Log\u00edstica do Transporte de Cargas. S\u00e3o Paulo: Atlas, 2001. 296 p.CAIXETA FILHO, J. V.; GAMEIRO, A. H. (org.). Sistemas de Gerenciamento de Transporte: Modelagem Matem\u00e1tica. S\u00e3o Paulo: Atlas, 2001. 125 p.CAIXETA FILHO, J. V. Pesquisa Operacional: T\u00e9cnicas de Otimiza\u00e7\u00e3o Aplicadas a Sistemas Agroindustriais. S\u00e3o Paulo: Atlas, 2001. 171 p.LEITE, P. R. Log\u00edstica Reversa: Competividade e Sustentabilidade. 3. ed. S\u00e3o Paulo: Saraiva, 2017. 360 p.\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "function Replace-Text($doc, $oldText, $newText) {\n    $range = $doc.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $result = $range.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    return $result\n}\n\nfunction Find-ParagraphIndex($doc, $text) {\n    $i = 0\n    foreach ($p in $doc.Paragraphs) {\n        $i = $i + 1\n        $t = $p.Range.Text.TrimEnd([char]13)\n        if ($t -eq $text) {\n            return $i\n        }\n    }\n    return -1\n}\n\nfunction Replace-ParagraphTextAndAddItalicAfter($doc, $oldText, $newText, $italicText) {\n    $idx = Find-ParagraphIndex $doc $oldText\n    $p = $doc.Paragraphs($idx)\n    $r = $p.Range\n    # Replace only the visible text, excluding the trailing paragraph mark,\n    # so the paragraph is not split in two.\n    $bodyOnly = $doc.Range($r.Start, $r.End - 1)\n    $bodyOnly.Text = $newText\n\n    # Insert a brand-new paragraph right after this one.\n    $r2 = $doc.Paragraphs($idx).Range\n    $r2.InsertParagraphAfter()\n    $newP = $doc.Paragraphs($idx + 1)\n    $newR = $newP.Range\n    $newR.Text = $italicText\n\n    # Apply italics only to the text run, not the paragraph mark.\n    $textOnlyRange = $doc.Range($newR.Start, $newR.End - 1)\n    $textOnlyRange.Font.Italic = 1\n}\n\n$d = $word.ActiveDocument\n\n# 1) Ativa\u00e7\u00e3o date update\nReplace-Text $d \"Ativa\u00e7\u00e3o: 01/01/2021\" \"Ativa\u00e7\u00e3o: 01/01/2024\" | Out-Null\n\n# 2) Objetivos: replace text + add italic English translation paragraph\nReplace-ParagraphTextAndAddItalicAfter $d `\n    \"Apresentar os conceitos de Log\u00edstica, Log\u00edstica Reversa e Gest\u00e3o da Cadeia de Suprimentos. Capacitar o aluno para aplica\u00e7\u00e3o de t\u00e9cnicas e m\u00e9todos quantitativos para otimiza\u00e7\u00e3o dos problemas em Log\u00edstica e Cadeias de Suprimentos.\" `\n    \"Fornecer conhecimentos que proporcionam uma vis\u00e3o hol\u00edstica a respeito da Cadeia de Suprimentos e da Log\u00edstica, apresentando m\u00e9todos e ferramentas para otimizar o desempenho das cadeias produtivas.\" `\n    \"Provide knowledge that paves a holistic view of Logistics and Supply Chain, presenting methods and tools available to optimize the performance of production chains.\"\n\n# 3) Docente respons\u00e1vel replacement\nReplace-Text $d \"5840560 - Marco Antonio Carvalho Pereira\" \"3295113 - Jos\u00e9 Eduardo Holler Branco\" | Out-Null\n\n# 4) Programa resumido: replace text + add italic English translation paragraph\nReplace-ParagraphTextAndAddItalicAfter $d `\n    \"1. Introdu\u00e7\u00e3o: 2. Gest\u00e3o estrat\u00e9gica3. Gest\u00e3o dos relacionamentos4. Gest\u00e3o global de suprimentos5. Avalia\u00e7\u00e3o de desempenho6. Mapeamento e an\u00e1lise de processos7. Gest\u00e3o de demanda8. Gest\u00e3o e coordena\u00e7\u00e3o de estoques9. Gest\u00e3o da log\u00edstica10. Log\u00edstica reversa\" `\n    \"Gerenciamento da cadeia de suprimentos e da log\u00edstica: planejamento, otimiza\u00e7\u00e3o e controle.\" `\n    \"Supply chain and logistics management: planning, optimization and control.\"\n\n# 5) Programa (full): replace text + add italic English translation paragraph\nReplace-ParagraphTextAndAddItalicAfter $d `\n    \"1. Introdu\u00e7\u00e3o: A concorr\u00eancia entre cadeias de suprimento. Defini\u00e7\u00e3o operacional. A globaliza\u00e7\u00e3o e a gest\u00e3o de cadeia de suprimentos. Governan\u00e7a das cadeias de suprimentos2. Gest\u00e3o estrat\u00e9gica: Estrat\u00e9gia de cadeia de suprimentos. Produtos funcionais x produtos inovadores. Fluxos empurrados puxados e h\u00edbridos. Custo de transa\u00e7\u00e3o e a decis\u00e3o estrat\u00e9gica de comprar ou fazer. Padroniza\u00e7\u00e3o. Integra\u00e7\u00e3o de parceiros da cadeia de suprimento no projeto de novos produtos e processos.3. Gest\u00e3o dos relacionamentos: Confian\u00e7a entre parceiros. Negocia\u00e7\u00e3o. Gest\u00e3o do relacionamento com clientes. Segmenta\u00e7\u00e3o de produtos. Gest\u00e3o do relacionamento com fornecedores4. Gest\u00e3o global de suprimentos: Tipos de suplemento. Estrutura organizacional para suprimentos. O processo de suprimento. Coopeti\u00e7\u00e3o. \u00c9tica e responsabilidade social na gest\u00e3o global de suprimentos5. Avalia\u00e7\u00e3o de desempenho: O que \u00e9 medi\u00e7\u00e3o de desempenho? Porque medir desempenho. Caracter\u00edsticas de uma boa medida de desempenho. Alinhamento de incentivos em cadeias globais de suprimento. Tipos de contrato de relacionamento6. Mapeamento e an\u00e1lise de processos: Principais processos na cadeia de suprimento. O modelo SCOR (Supply Chain Operations Reference). An\u00e1lise e melhoramento de processos.7. Gest\u00e3o de demanda: A\u00e7\u00f5es sobre a demanda para redu\u00e7\u00e3o de variabilidade. Causas da variabilidade da demanda. Previs\u00e3o de demanda. Processo de previs\u00e3o de vendas. M\u00e9todos usados em previs\u00f5es. M\u00e9todo Delphi. Incerteza de previs\u00e3o8. Gest\u00e3o e coordena\u00e7\u00e3o de estoques: Defini\u00e7\u00e3o de estoques. Causa do surgimento dos estoques. Tipos de estoque. VMI (vendor management inventory) - estoque gerenciado pelo distribuidor. VOI (vendor owner inventory) - consigna\u00e7\u00e3o9. Gest\u00e3o da log\u00edstica: Centraliza\u00e7\u00e3o versus descentraliza\u00e7\u00e3o. Pontos de armazenagem/distribui\u00e7\u00e3o. Fun\u00e7\u00f5es dos armaz\u00e9ns. Sistemas log\u00edsticos escalonados. Localiza\u00e7\u00e3o de unidades log\u00edsticas. Gest\u00e3o de transportes na cadeia de suprimentos.10. Log\u00edstica reversa: Conceito, import\u00e2ncia, estrutura e tend\u00eancias. Sustentabilidade. Ciclo fechado. Tipos de ciclo fechado. Motiva\u00e7\u00e3o empresarial. Gerenciamento integrado de res\u00edduos. Modelos de roteiriza\u00e7\u00e3o. Programa\u00e7\u00e3o de frotas de ve\u00edculos.\" `\n    \"i) Introdu\u00e7\u00e3o \u00e0 Log\u00edstica e Cadeia de Suprimentos; ii) Planejamento da cadeia de suprimentos; iii) Planejamento do transporte; iv) Custos log\u00edsticos; v) Tipos de cargas e sistemas de armazenamento; vi) Modelos de transporte; vii) Modelos de localiza\u00e7\u00e3o; viii) Planejamento do estoque; ix) Log\u00edstica Reversa e Economia Circular; e x)  Controle da log\u00edstica e cadeia de suprimentos.\" `\n    \"i) Introduction to Logistics and Supply Chain; ii) Supply chain planning; iii) Transport planning; iv) Logistic costs; v) Types of cargos and storage systems; vi) Stock planning; vii) Transport models; viii) Location models; ix) Reverse Logistics and Circular Economy; and x) Control of logistics and supply chain.\"\n\n# 6) M\u00e9todo: text replacement\nReplace-Text $d \"Aulas expositivas te\u00f3ricas, aulas de exerc\u00edcios.\" \"Provas, trabalhos em grupo, exerc\u00edcios individuais e semin\u00e1rios.\" | Out-Null\n\n# 7) Crit\u00e9rio: text replacement\nReplace-Text $d \"M\u00e9dia de Provas e trabalhos (MF).\" \"M\u00e9dia das atividades avaliativas.\" | Out-Null\n\n# 8) Norma de recupera\u00e7\u00e3o: text replacement\nReplace-Text $d \"Prova de Recupera\u00e7\u00e3o (PR). A Nota final (NF) ser\u00e1 a m\u00e9dia aritm\u00e9tica entre MF e PR\" \"NF = (MF + PR)/2, onde MF \u00e9 a m\u00e9dia final da avalia\u00e7\u00e3o e PR \u00e9 uma prova de recupera\u00e7\u00e3o.\" | Out-Null\n\n# 9) Bibliografia: text replacement\nReplace-Text $d \"CORR\u00caA, HENRIQUE LUIZ. Gest\u00e3o de rede de suprimento: integrando cadeias de suprimento no mundo globalizado. Editora Atlas, 2009CORREA, HENRIQUE LUIZ. Administra\u00e7\u00e3o de cadeias de suprimento e log\u00edstica: o essencial. Editora Atlas 2014PIRES, S\u00c9RGIO. Gest\u00e3o da cadeia de suprimentos (Supply Chain Management): conceitos, estrat\u00e9gias, pr\u00e1ticas e casos. Editora Atlas segunda edi\u00e7\u00e3o. 2009IYER, ANANTH; SESHHADRI, SHIDHAR; VASHER, ROY. A gest\u00e3o da cadeia de suprimentos da Toyota. Bookman. 2009MYERSON, PAUL. Lean Supply Chain and logistics management. McGrawHill. 2012\" \"BOWERSOX, D. J.; CLOSS, D. J.; COOPER, M. B.; BOWERSOX, J. C. Gest\u00e3o Log\u00edstica da Cadeia de Suprimentos. 4. ed. AMGH, 2013. 472 p.BARTHOLOMEU, D. B.; CAIXETA FILHO, J. V. Log\u00edstica Ambiental de Res\u00edduos S\u00f3lidos. S\u00e3o Paulo: Atlas, 2011, 249 p.CHOPRA, S.; MEINDL, P. Gest\u00e3o da cadeia de suprimentos: estrat\u00e9gia, planejamento e Opera\u00e7\u00f5es. 6. ed. Pearson, 2015. 544 p.CAIXETA FILHO, J. V.; MARTINS, R. S. (org.). Gest\u00e3o Log\u00edstica do Transporte de Cargas. S\u00e3o Paulo: Atlas, 2001. 296 p.CAIXETA FILHO, J. V.; GAMEIRO, A. H. (org.). Sistemas de Gerenciamento de Transporte: Modelagem Matem\u00e1tica. S\u00e3o Paulo: Atlas, 2001. 125 p.CAIXETA FILHO, J. V. Pesquisa Operacional: T\u00e9cnicas de Otimiza\u00e7\u00e3o Aplicadas a Sistemas Agroindustriais. S\u00e3o Paulo: Atlas, 2001. 171 p.LEITE, P. R. Log\u00edstica Reversa: Competividade e Sustentabilidade. 3. ed. S\u00e3o Paulo: Saraiva, 2017. 360 p.\" | Out-Null\n"}
